$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/3/2025  Through  11/9/2025"

# --- Data cell updates ---
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = -25
$ws.Range("L14").Value = -25
$ws.Range("M14").Value = -81.25
$ws.Range("N14").Value = -88
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("H15").Value = -100
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 55.555555555555
$ws.Range("L15").Value = -17.647058823529
$ws.Range("M15").Value = -36.363636363636
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("J14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = -75
$ws.Range("K14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 102
$ws.Range("J16").Value = 112
$ws.Range("K16").Value = -8.928571428571
$ws.Range("L16").Value = -21.538461538461
$ws.Range("M16").Value = -65.306122448979
$ws.Range("N16").Value = -90.431519699812
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -13.043478260869
$ws.Range("I17").Value = 207
$ws.Range("J17").Value = 260
$ws.Range("K17").Value = -20.384615384615
$ws.Range("L17").Value = -20.384615384615
$ws.Range("M17").Value = -26.071428571428
$ws.Range("N17").Value = -72.619047619047
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 6
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 79
$ws.Range("J18").Value = 71
$ws.Range("K18").Value = 11.267605633802
$ws.Range("L18").Value = -34.710743801652
$ws.Range("M18").Value = -61.83574879227
$ws.Range("N18").Value = -88.665710186513
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = -11.111111111111
$ws.Range("I19").Value = 239
$ws.Range("J19").Value = 238
$ws.Range("K19").Value = 0.420168067226
$ws.Range("L19").Value = -21.381578947368
$ws.Range("M19").Value = -20.333333333333
$ws.Range("N19").Value = -25.3125
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 75
$ws.Range("J20").Value = 86
$ws.Range("K20").Value = -12.790697674418
$ws.Range("L20").Value = -20.212765957446
$ws.Range("M20").Value = 1.351351351351
$ws.Range("N20").Value = -85.65965583174
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 60
$ws.Range("H21").Value = -5
$ws.Range("I21").Value = 719
$ws.Range("J21").Value = 780
$ws.Range("K21").Value = -7.820512820512
$ws.Range("L21").Value = -22.68817204301
$ws.Range("M21").Value = -39.731768650461
$ws.Range("N21").Value = -79.201619901648
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("H22").Value = -100
$ws.Range("M22").Value = -39.130434782608
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 120
$ws.Range("I23").Value = 95
$ws.Range("J23").Value = 83
$ws.Range("K23").Value = 14.457831325301
$ws.Range("L23").Value = 25
$ws.Range("M23").Value = 33.802816901408
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -5.882352941176
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 62
$ws.Range("H24").Value = -11.290322580645
$ws.Range("I24").Value = 620
$ws.Range("J24").Value = 621
$ws.Range("K24").Value = -0.161030595813
$ws.Range("L24").Value = -16.554508748317
$ws.Range("M24").Value = -11.80654338549
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 12
$ws.Range("H25").Value = -14.285714285714
$ws.Range("I25").Value = 72
$ws.Range("J25").Value = 119
$ws.Range("K25").Value = -39.495798319327
$ws.Range("L25").Value = -59.090909090909
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = 11.904761904761
$ws.Range("I26").Value = 386
$ws.Range("J26").Value = 368
$ws.Range("K26").Value = 4.891304347826
$ws.Range("L26").Value = -10.854503464203
$ws.Range("M26").Value = -46.08938547486
$ws.Range("C27").Value = 1
$ws.Range("J14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("I27").Value = 20
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = -20
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 32
$ws.Range("J28").Value = 33
$ws.Range("K28").Value = -3.030303030303
$ws.Range("L28").Value = 14.285714285714
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("N29").Value = -90.184049079754
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("N30").Value = -91.836734693877
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"
$ws.Range("G31").Copy()
$ws.Range("F31").PasteSpecial(-4122)

# --- Column width cosmetic tweak (bestFit narrowed after data refresh) ---
$ws.Columns.Item(5).ColumnWidth = 6.142857142857143
$ws.Columns.Item(9).ColumnWidth = 6.142857142857143
